$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 'Duel Decks Ajani vs. Nicol Bolas Tokens (TDDH)'
$ws.Range("A2").Value = 'Aethersnipe'
$ws.Range("A3").Value = 'Air Elemental'
$ws.Range("A4").Value = 'Ancestral Vision'
$ws.Range("A5").Value = 'Bottle Gnomes'
$ws.Range("A6").Value = 'Brine Elemental'
$ws.Range("A7").Value = 'Chandra Nalaar'
$ws.Range("A8").Value = 'Chartooth Cougar'
$ws.Range("A9").Value = 'Condescend'
$ws.Range("A10").Value = 'Cone of Flame'
$ws.Range("A11").Value = 'Counterspell'
$ws.Range("A12").Value = 'Daze'
$ws.Range("A13").Value = 'Demonfire'
$ws.Range("A14").Value = 'Errant Ephemeron'
$ws.Range("A15").Value = 'Fact or Fiction'
$ws.Range("A16").Value = 'Fathom Seer'
$ws.Range("A17").Value = 'Fireball'
$ws.Range("A18").Value = 'Fireblast'
$ws.Range("A19").Value = 'Firebolt'
$ws.Range("A20").Value = 'Fireslinger'
$ws.Range("A21").Value = 'Flame Javelin'
$ws.Range("A22").Value = 'Flamekin Brawler'
$ws.Range("A23").Value = 'Flametongue Kavu'
$ws.Range("A24").Value = 'Flamewave Invoker'
$ws.Range("A25").Value = 'Fledgling Mawcor'
$ws.Range("A26").Value = 'Furnace Whelp'
$ws.Range("A27").Value = 'Guile'
$ws.Range("A28").Value = 'Gush'
$ws.Range("A29").Value = 'Hostility'
$ws.Range("A30").Value = 'Incinerate'
$ws.Range("A31").Value = 'Ingot Chewer'
$ws.Range("A32").Value = 'Inner-Flame Acolyte'
$ws.Range("A33").Value = 'Island'
$ws.Range("A34").Value = 'Island'
$ws.Range("A35").Value = 'Island'
$ws.Range("A36").Value = 'Island'
$ws.Range("A37").Value = 'Jace Beleren'
$ws.Range("A38").Value = 'Keldon Megaliths'
$ws.Range("A39").Value = 'Magma Jet'
$ws.Range("A40").Value = 'Man-o''-War'
$ws.Range("A41").Value = 'Martyr of Frost'
$ws.Range("A42").Value = 'Mind Stone'
$ws.Range("A43").Value = 'Mountain'
$ws.Range("A44").Value = 'Mountain'
$ws.Range("A45").Value = 'Mountain'
$ws.Range("A46").Value = 'Mountain'
$ws.Range("A47").Value = 'Mulldrifter'
$ws.Range("A48").Value = 'Ophidian'
$ws.Range("A49").Value = 'Oxidda Golem'
$ws.Range("A50").Value = 'Pyre Charger'
$ws.Range("A51").Value = 'Quicksilver Dragon'
$ws.Range("A52").Value = 'Rakdos Pit Dragon'
$ws.Range("A53").Value = 'Repulse'
$ws.Range("A54").Value = 'Riftwing Cloudskate'
$ws.Range("A55").Value = 'Seal of Fire'
$ws.Range("A56").Value = 'Slith Firewalker'
$ws.Range("A57").Value = 'Soulbright Flamekin'
$ws.Range("A58").Value = 'Spire Golem'
$ws.Range("A59").Value = 'Terrain Generator'
$ws.Range("A60").Value = 'Voidmage Apprentice'
$ws.Range("A61").Value = 'Wall of Deceit'
$ws.Range("A62").Value = 'Waterspout Djinn'
$ws.Range("A63").Value = 'Willbender'

$ws.Range("A64").EntireRow.Delete()
